$d = $word.ActiveDocument

# Locate the paragraph that contains the "LOQ4086" requirement line; the
# three paragraphs that directly follow it (a blank paragraph, the "Ver no
# Jupiter..." paragraph, and the "(c) 2020 ..." footer paragraph) must be
# removed, while the page-break paragraph further down stays untouched.
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*LOQ4086*") {
        $anchor = $i
        break
    }
}

if ($anchor -ne $null) {
    # Delete the three trailing paragraphs from the bottom up so that the
    # indices of the ones not yet removed stay valid.
    $d.Paragraphs.Item($anchor + 3).Range.Delete()   # "(c) 2020 ... Creative Commons Attribution"
    $d.Paragraphs.Item($anchor + 2).Range.Delete()   # "Ver no Jupiter Salvar em pdf Salvar em docx"
    $d.Paragraphs.Item($anchor + 1).Range.Delete()   # blank paragraph right after LOQ4086 line
}
